$wb = $excel.ActiveWorkbook

# The "States" worksheet (5th tab) becomes the active sheet/tab, with cell
# E7 selected as the active cell.
$ws = $wb.Worksheets.Item(5)
$ws.Activate()
$ws.Range("E7").Select()

# Enter the attrition-rate percentages for the first five state rows.
# These cells pick up the workbook's 0.00% number format (style 19) when
# the format/value moves in from the previously plain "0" (style 17).
$ws.Range("E2").NumberFormat = "0.00%"
$ws.Range("E2").Value = 0.02

$ws.Range("E3").NumberFormat = "0.00%"
$ws.Range("E3").Value = 0.03

$ws.Range("E4").NumberFormat = "0.00%"
$ws.Range("E4").Value = 0.02

$ws.Range("E5").Value = 0.01

$ws.Range("E6").Value = 0.01
